# The deck ships with two themes: the slide master (ppt/theme/theme1.xml)
# carries the "Integral" / "Red Violet" palette while the notes master
# (ppt/theme/theme2.xml) carries the stock "Office Theme" palette. The
# commit swaps them: the slide master is re-themed to the default
# "Office Theme" colours.
#
# PowerPoint's object model exposes the twelve theme colour slots via
# SlideMaster.Theme.ThemeColorScheme(1..12), in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB is stored as the usual OLE BGR-packed long (same encoding the
# VBA RGB() function produces), so each hex colour below is expressed
# as b*65536 + g*256 + r.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}

# These renames mirror what a full theme-gallery swap would also change
# (theme name "Integral" -> "Office Theme", colour-scheme name
# "Red Violet" -> "Office"); harmless to attempt even where read-only.
$tcs.Name = "Office"
$p.Designs.Item(1).Name = "Office Theme"
